$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2697.8
$ws.Range("I58").Value = 526.3333
$ws.Range("J58").Value = 3628.4285
$ws.Range("K58").Value = 1578.9999
$ws.Range("L58").Value = 10885.2855
$ws.Range("M58").Value = -1428.9999
$ws.Range("N58").Value = -11185.2855

$ws.Range("H69").Value = 4713.28
$ws.Range("I69").Value = 4092.4614
$ws.Range("J69").Value = 5385.8335
$ws.Range("K69").Value = 12277.3842
$ws.Range("L69").Value = 16157.5005
$ws.Range("M69").Value = -11403.3842
$ws.Range("N69").Value = -17905.5005

$ws.Range("H72").Value = 4713.28
$ws.Range("I72").Value = 4092.4614
$ws.Range("J72").Value = 5385.8335
$ws.Range("K72").Value = 36832.1526
$ws.Range("L72").Value = 48472.5015
$ws.Range("M72").Value = -32464.1526
$ws.Range("N72").Value = -57208.5015

$ws.Range("H76").Value = 4525.9165
$ws.Range("I76").Value = 3340.6
$ws.Range("J76").Value = 5372.5713
$ws.Range("K76").Value = 3340.6
$ws.Range("L76").Value = 5372.5713
$ws.Range("M76").Value = -3025.6
$ws.Range("N76").Value = -6002.5713

$ws.Range("H79").Value = 4525.9165
$ws.Range("I79").Value = 3340.6
$ws.Range("J79").Value = 5372.5713
$ws.Range("K79").Value = 3340.6
$ws.Range("L79").Value = 5372.5713
$ws.Range("M79").Value = -2248.6
$ws.Range("N79").Value = -7556.5713

$ws.Range("H100").Value = 6411651.5
$ws.Range("I100").Value = 15152670
$ws.Range("J100").Value = 1571.0667
$ws.Range("K100").Value = 15152670
$ws.Range("L100").Value = 1571.0667
$ws.Range("M100").Value = -15152129
$ws.Range("N100").Value = -2653.0667

$ws.Range("H106").Value = 333338340
$ws.Range("I106").Value = 111117784
$ws.Range("K106").Value = 111117784
$ws.Range("M106").Value = -111117153

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H112").Value = 7813525
$ws.Range("J112").Value = 8475630
$ws.Range("L112").Value = 25426890
$ws.Range("N112").Value = -25429106

$ws.Range("H129").Value = 927.54
$ws.Range("J129").Value = 1006.375
$ws.Range("L129").Value = 3019.125
$ws.Range("N129").Value = -13019.125

$ws.Range("H132").Value = 1332.2174
$ws.Range("I132").Value = 1347.3636
$ws.Range("K132").Value = 4042.0908
$ws.Range("M132").Value = -1512.0908

$ws.Range("H135").Value = 2067.625
$ws.Range("I135").Value = 2138.8
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 19249.2
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -16714.2
$ws.Range("N135").Value = -14070

$ws.Range("H138").Value = 5579.1665
$ws.Range("I138").Value = 1540.4286
$ws.Range("J138").Value = 7242.1763
$ws.Range("K138").Value = 4621.2858
$ws.Range("L138").Value = 21726.5289
$ws.Range("M138").Value = 518.7142000000003
$ws.Range("N138").Value = -32006.5289

$ws.Range("H141").Value = 2930
$ws.Range("I141").Value = 2410
$ws.Range("J141").Value = 4750
$ws.Range("K141").Value = 7230
$ws.Range("L141").Value = 14250
$ws.Range("M141").Value = -2050
$ws.Range("N141").Value = -24610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10139.826
$ws.Range("I45").Value = 20802.7
$ws.Range("J45").Value = 1937.6154
$ws.Range("K45").Value = 20802.7
$ws.Range("L45").Value = 1937.6154
$ws.Range("M45").Value = -20425.7
$ws.Range("N45").Value = -2691.6154

$ws.Range("H110").Value = 831.9375
$ws.Range("I110").Value = 809.25
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 809.25
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 1235.75
$ws.Range("N110").Value = -4990

$ws.Range("H132").Value = 6883
$ws.Range("I132").Value = 1492.7858
$ws.Range("J132").Value = 11913.866
$ws.Range("K132").Value = 4478.357400000001
$ws.Range("L132").Value = 35741.598
$ws.Range("M132").Value = -1948.357400000001
$ws.Range("N132").Value = -40801.598

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20835190
$ws.Range("I86").Value = 27779628
$ws.Range("J86").Value = 1872.25
$ws.Range("K86").Value = 27779628
$ws.Range("L86").Value = 1872.25
$ws.Range("M86").Value = -27778505
$ws.Range("N86").Value = -4118.25

$ws.Range("H89").Value = 20835190
$ws.Range("I89").Value = 27779628
$ws.Range("J89").Value = 1872.25
$ws.Range("K89").Value = 138898140
$ws.Range("L89").Value = 9361.25
$ws.Range("M89").Value = -138892524
$ws.Range("N89").Value = -20593.25

$ws.Range("H107").Value = 1006.1739
$ws.Range("I107").Value = 1019.0625
$ws.Range("K107").Value = 1019.0625
$ws.Range("M107").Value = 900.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 413.7143
$ws.Range("I7").Value = 449.5
$ws.Range("J7").Value = 199
$ws.Range("K7").Value = 449.5
$ws.Range("L7").Value = 199
$ws.Range("M7").Value = -336.5
$ws.Range("N7").Value = -425

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314

$ws.Range("H122").Value = 1098.75
$ws.Range("I122").Value = 1263.3334
$ws.Range("K122").Value = 3790.0002
$ws.Range("M122").Value = -1340.0002

$ws.Range("H134").Value = 1569.6364
$ws.Range("I134").Value = 1501.6842
$ws.Range("K134").Value = 4505.0526
$ws.Range("M134").Value = -1970.0526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1591.25
$ws.Range("J13").Value = 2997.5
$ws.Range("L13").Value = 8992.5
$ws.Range("N13").Value = -9328.5

$ws.Range("H31").Value = 625
$ws.Range("I31").Value = 500.33334
$ws.Range("J31").Value = 999
$ws.Range("K31").Value = 1501.00002
$ws.Range("L31").Value = 2997
$ws.Range("M31").Value = -1213.00002
$ws.Range("N31").Value = -3573

$ws.Range("H63").Value = 4800
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 5400
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 16200
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -17698

$ws.Range("H66").Value = 4800
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 5400
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 48600
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -56088

$ws.Range("H129").Value = 16668095
$ws.Range("J129").Value = 3160
$ws.Range("L129").Value = 9480
$ws.Range("N129").Value = -19480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3459
$ws.Range("J80").Value = 3100
$ws.Range("L80").Value = 3100
$ws.Range("N80").Value = -5096

$ws.Range("H83").Value = 3459
$ws.Range("J83").Value = 3100
$ws.Range("L83").Value = 15500
$ws.Range("N83").Value = -25484

$ws.Range("H102").Value = 1490.5454
$ws.Range("I102").Value = 1573.375
$ws.Range("J102").Value = 1269.6666
$ws.Range("K102").Value = 1573.375
$ws.Range("L102").Value = 1269.6666
$ws.Range("M102").Value = 48.625
$ws.Range("N102").Value = -4513.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2925.3333
$ws.Range("I7").Value = 1699.8572
$ws.Range("K7").Value = 1699.8572
$ws.Range("M7").Value = -1587.8572

$ws.Range("H22").Value = 3970749.5
$ws.Range("I22").Value = 37037372
$ws.Range("J22").Value = 2755.08
$ws.Range("K22").Value = 37037372
$ws.Range("L22").Value = 2755.08
$ws.Range("M22").Value = -37037077
$ws.Range("N22").Value = -3345.08

$ws.Range("H27").Value = 3970749.5
$ws.Range("I27").Value = 37037372
$ws.Range("J27").Value = 2755.08
$ws.Range("K27").Value = 37037372
$ws.Range("L27").Value = 2755.08
$ws.Range("M27").Value = -37037265
$ws.Range("N27").Value = -2969.08

$ws.Range("H40").Value = 200003980
$ws.Range("I40").Value = 500002500
$ws.Range("K40").Value = 500002500
$ws.Range("M40").Value = -500002364

$ws.Range("H46").Value = 23810654
$ws.Range("I46").Value = 47619944
$ws.Range("J46").Value = 1364.1428
$ws.Range("K46").Value = 47619944
$ws.Range("L46").Value = 1364.1428
$ws.Range("M46").Value = -47619756
$ws.Range("N46").Value = -1740.1428

$ws.Range("H55").Value = 65217740
$ws.Range("I55").Value = 83333520
$ws.Range("J55").Value = 45455070
$ws.Range("K55").Value = 83333520
$ws.Range("L55").Value = 45455070
$ws.Range("M55").Value = -83333347
$ws.Range("N55").Value = -45455416

$ws.Range("H93").Value = 41684740
$ws.Range("I93").Value = 23266.666
$ws.Range("J93").Value = 166669170
$ws.Range("K93").Value = 23266.666
$ws.Range("L93").Value = 166669170
$ws.Range("M93").Value = -22018.666
$ws.Range("N93").Value = -166671666

$ws.Range("H122").Value = 5497359.5
$ws.Range("I122").Value = 10206954
$ws.Range("J122").Value = 2832.5
$ws.Range("K122").Value = 30620862
$ws.Range("L122").Value = 8497.5
$ws.Range("M122").Value = -30618412
$ws.Range("N122").Value = -13397.5

$ws.Range("H126").Value = 2925.3333
$ws.Range("I126").Value = 1699.8572
$ws.Range("K126").Value = 5099.571599999999
$ws.Range("M126").Value = -2629.571599999999

$ws.Range("H132").Value = 16055822
$ws.Range("I132").Value = 22814826
$ws.Range("J132").Value = 3187.875
$ws.Range("K132").Value = 68444478
$ws.Range("L132").Value = 9563.625
$ws.Range("M132").Value = -68441948
$ws.Range("N132").Value = -14623.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
